$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Workbook window view (best-effort; engine may not persist xWindow/yWindow) ---
$wb.Windows.Item(1).Left = 4880
$wb.Windows.Item(1).Top = 980

# --- Update data values and labels on the "definition" sheet ---
# Row 1
$ws.Range("A1").Value = "Condenser"
$ws.Range("B1").Value = "definition"

# Row 2
$ws.Range("A2").Value = 201
$ws.Range("B2").Value = "201 = user defined, 202 = circuited tube pattern, 203 = single finned tube, "

# Row 3
$ws.Range("A3").Value = 0.0075
$ws.Range("B3").Value = "tube inside diameter (m)"

# Row 4
$ws.Range("A4").Value = 0.61
$ws.Range("B4").Value = "tube length (m)"

# Row 5
$ws.Range("A5").Value = 0.0002
$ws.Range("B5").Value = "wall thickness of pipe (m), for micro-fin tubes, it doesn't include the fin height and the wall thickness"

# Row 6 (style changes to an explicit black-colored font)
$ws.Range("A6").Value = 0.0287
$ws.Range("A6").Font.Color = 0
$ws.Range("B6").Value = "outside fin diameter (m)"

# Row 7
$ws.Range("A7").Value = 0.0015
$ws.Range("B7").Value = "fin pitch (m), fin space+one fin thickness"

# Row 8
$ws.Range("A8").Value = 0.00013
$ws.Range("B8").Value = "fin thickness (m)"

# Row 10 label is set first so that the new shared string for "(P_l)" is
# created before the new shared string for "(P_t)" (matches target ordering
# of the shared strings table).
$ws.Range("B10").Value = "spacing between tubs in the longitudual direction (m) (P_l)"

# Row 9
$ws.Range("A9").Value = 0.024211
$ws.Range("B9").Value = "spacing between tubes in bank (m) (P_t)"

# Row 10 (continued)
$ws.Range("A10").Value = 0.0125

# Row 11
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "number of segments per tube in finite difference model"

# Row 12
$ws.Range("A12").Value = 0.0127
$ws.Range("B12").Value = "radius of return bend (m)"

# Row 13
$ws.Range("A13").Value = -1
$ws.Range("B13").Value = "nominal air mass flux (kg/s/m^2), correponding to the maximum air flux"

# Row 14
$ws.Range("A14").Value = 1
$ws.Range("B14").Value = "number of parallel branches in main section (circuit pattern model only)"

# Row 15
$ws.Range("A15").Value = 54
$ws.Range("B15").Value = "number of tubes in main branch (circuit pattern model only)"

# Row 16
$ws.Range("A16").Value = 0
$ws.Range("B16").Value = "number of tubes in subcooling branch (circuit pattern model only)"

# Row 17
$ws.Range("A17").Value = 3
$ws.Range("D17").Value = "number of banks"

# Row 18
$ws.Range("A18").Value = 0.281
$ws.Range("D18").Value = "Frontal area [m2]"

# Row 19
$ws.Range("A19").Value = 1
$ws.Range("D19").Value = "air side convection coefficient, for the first series test forTxV it is 1.0, since some unaccounted operation deviation"

# Row 20
$ws.Range("A20").Value = 1
$ws.Range("C20").Value = "refrigerant side two-phase flow convection coefficient"

# Row 21
$ws.Range("A21").Value = 1
$ws.Range("C21").Value = "refrigerant side subcooling convection coefficient"

# Row 22
$ws.Range("A22").Value = 1
$ws.Range("C22").Value = "refrigerant side pressure drop coefficient"

# Row 23
$ws.Range("A23").Value = 1
$ws.Range("C23").Value = "fouling factor (1.0-no, 0.0-completely)"

# Row 24
$ws.Range("A24").Value = 0
$ws.Range("D24").Value = "microfin type, 0=smooth tube, 1=helical, 2=cross-grooved, 3=herringbone"

# Row 25
$ws.Range("A25").Value = 0
$ws.Range("D25").Value = "fin number in a micro-fin tube"

# Row 26
$ws.Range("A26").Value = 0
$ws.Range("D26").Value = "fin apex angle in a micro-fin tube"

# Row 27
$ws.Range("A27").Value = 0
$ws.Range("D27").Value = "fin helix angle in a micro-fin tube"

# Row 28
$ws.Range("A28").Value = 0
$ws.Range("C28").Value = "fin height in a micro-fin tube"

# Row 29
$ws.Range("A29").Value = 0
$ws.Range("C29").Value = "base width of a single fin"

# Row 30
$ws.Range("A30").Value = 0
$ws.Range("D30").Value = "top width of a single fin"

# Row 31
$ws.Range("A31").Value = 0
$ws.Range("D31").Value = "base distance between two neighboring fins"

# Row 32
$ws.Range("A32").Value = 5
$ws.Range("D32").Value = "airside fin type, 1-plain, 2-corrugated, 3-slit, 4-louvered, 5-convex louvered, 6-smooth wavy, 7-spine"

# Row 33
$ws.Range("A33").Value = 0
$ws.Range("D33").Value = "substructure of fin surface"

# Row 34
$ws.Range("A34").Value = 0
$ws.Range("D34").Value = "substructure of fin surface"

# Row 35
$ws.Range("A35").Value = 0
$ws.Range("D35").Value = "substructure of fin surface"

# Row 36
$ws.Range("A36").Value = 0
$ws.Range("D36").Value = "substructure of fin surface"

# Row 37
$ws.Range("A37").Value = 0
$ws.Range("D37").Value = "substructure of fin surface"

# Row 38
$ws.Range("A38").Value = 400
$ws.Range("D38").Value = "tube wall conductance"

# Row 39
$ws.Range("A39").Value = 237
$ws.Range("D39").Value = "fin conductance"

# --- Selection as shown in final sheetView ---
$ws.Range("A11").Select()
